$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Cell value / formula text updates
# ---------------------------------------------------------------------------

# Row 2 used to be the "CasesTab" row; it is now "ParticipantsTab"
$ws.Range("A2").Value = "ParticipantsTab"

$participantsQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE p.gender in ['Unknown']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct coalesce(samp.sample_id, "Not specified in data"))) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@
$ws.Range("B2").Value = $participantsQuery

$statQuery = @'
CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE p.gender in ['Unknown']
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE p.gender in ['Unknown']
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE p.gender in ['Unknown']
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`
'@
$ws.Range("C2").Value = $statQuery

# Row 3 ("SamplesTab") - detail query text updated, tab name unchanged
$samplesQuery = @'
MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE p.gender in ["Unknown"]
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(samp.sample_tumor_status,'') as `Tumor`,
    coalesce(samp.sample_type,'') as `Analyte Type`
ORDER BY samp.sample_id limit 100
'@
$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statQuery

# Row 4 ("FilesTab") - detail query text updated, tab name unchanged
$filesQuery = @'
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE p.gender in ['Unknown']
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name,'') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id, '') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER BY f.file_name limit 100
'@
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery

# ---------------------------------------------------------------------------
# 2) Font size: whole used range moves from 11/12pt to 14pt
# ---------------------------------------------------------------------------

$noWrapCells = @("A1","B1","C1","D1","E1","A2","D2","E2","A3","D3","E3","A4","D4","E4")
foreach ($c in $noWrapCells) {
    $ws.Range($c).Font.Size = 14
}

$wrapCells = @("B2","C2","B3","C3","B4","C4","B5","C5","C6")
foreach ($c in $wrapCells) {
    $ws.Range($c).Font.Size = 14
    $ws.Range($c).WrapText = $true
}

# ---------------------------------------------------------------------------
# 3) Row heights for the (now much longer) query rows
# ---------------------------------------------------------------------------

$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5

# ---------------------------------------------------------------------------
# 4) Column widths (best effort - headless engine quantizes to char units)
# ---------------------------------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 19.57
$ws.Columns.Item(2).ColumnWidth = 75.86
$ws.Columns.Item(3).ColumnWidth = 57.86
$ws.Columns.Item(4).ColumnWidth = 60
$ws.Columns.Item(5).ColumnWidth = 58.57

# ---------------------------------------------------------------------------
# 5) Selection moves from B2 to D3
# ---------------------------------------------------------------------------

$ws.Range("D3").Select() | Out-Null

Write-Host "edit complete"
